# Update market price / profit data cells per scheduled data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 125000610
$ws.Range("I33").Value = 125000610
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 125000610
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -125000381
$ws.Range("N33").ClearContents()

$ws.Range("H40").Value = 1868.9474
$ws.Range("I40").Value = 1772.2222
$ws.Range("J40").Value = 1956
$ws.Range("K40").Value = 1772.2222
$ws.Range("L40").Value = 1956
$ws.Range("M40").Value = -1597.2222
$ws.Range("N40").Value = -2306

$ws.Range("H64").Value = 3455.1428
$ws.Range("J64").Value = 3426.6
$ws.Range("L64").Value = 3426.6
$ws.Range("N64").Value = -3922.6

$ws.Range("H67").Value = 3455.1428
$ws.Range("J67").Value = 3426.6
$ws.Range("L67").Value = 3426.6
$ws.Range("N67").Value = -5142.6

$ws.Range("H137").Value = 1323.4073
$ws.Range("I137").Value = 970.4737
$ws.Range("K137").Value = 2911.4211
$ws.Range("M137").Value = -361.4211

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4121.6875
$ws.Range("I61").Value = 4196.467
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 4196.467
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -3984.467
$ws.Range("N61").Value = -3424

$ws.Range("H74").Value = 4757.05
$ws.Range("I74").Value = 3091.4
$ws.Range("J74").Value = 6422.7
$ws.Range("K74").Value = 3091.4
$ws.Range("L74").Value = 6422.7
$ws.Range("M74").Value = -2217.4
$ws.Range("N74").Value = -8170.7

$ws.Range("H77").Value = 4757.05
$ws.Range("I77").Value = 3091.4
$ws.Range("J77").Value = 6422.7
$ws.Range("K77").Value = 15457
$ws.Range("L77").Value = 32113.5
$ws.Range("M77").Value = -11089
$ws.Range("N77").Value = -40849.5

$ws.Range("H132").Value = 14709414
$ws.Range("I132").Value = 16669668
$ws.Range("J132").Value = 7506.5
$ws.Range("K132").Value = 50009004
$ws.Range("L132").Value = 22519.5
$ws.Range("M132").Value = -50006474
$ws.Range("N132").Value = -27579.5

$ws.Range("H136").Value = 4121.6875
$ws.Range("I136").Value = 4196.467
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 12589.401
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -10039.401
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 821.2
$ws.Range("I7").Value = 821.2
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 821.2
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -708.2
$ws.Range("N7").ClearContents()

$ws.Range("H75").Value = 45118
$ws.Range("I75").Value = 40000
$ws.Range("J75").Value = 45849.145
$ws.Range("K75").Value = 40000
$ws.Range("L75").Value = 45849.145
$ws.Range("M75").Value = -39064
$ws.Range("N75").Value = -47721.145

$ws.Range("H78").Value = 45118
$ws.Range("I78").Value = 40000
$ws.Range("J78").Value = 45849.145
$ws.Range("K78").Value = 120000
$ws.Range("L78").Value = 137547.435
$ws.Range("M78").Value = -115320
$ws.Range("N78").Value = -146907.435

$ws.Range("H105").Value = 4344.294
$ws.Range("I105").Value = 3250.875
$ws.Range("J105").Value = 4680.731
$ws.Range("K105").Value = 3250.875
$ws.Range("L105").Value = 4680.731
$ws.Range("M105").Value = -1503.875
$ws.Range("N105").Value = -8174.731

$ws.Range("H134").Value = 2403.5854
$ws.Range("I134").Value = 1459.6389
$ws.Range("K134").Value = 4378.9167
$ws.Range("M134").Value = -1843.9167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4605.971
$ws.Range("I31").Value = 5944.04
$ws.Range("J31").Value = 1260.8
$ws.Range("K31").Value = 5944.04
$ws.Range("L31").Value = 1260.8
$ws.Range("M31").Value = -5649.04
$ws.Range("N31").Value = -1850.8

$ws.Range("H34").Value = 4605.971
$ws.Range("I34").Value = 5944.04
$ws.Range("J34").Value = 1260.8
$ws.Range("K34").Value = 5944.04
$ws.Range("L34").Value = 1260.8
$ws.Range("M34").Value = -5742.04
$ws.Range("N34").Value = -1664.8

$ws.Range("H58").Value = 2084.6316
$ws.Range("I58").Value = 650.5714
$ws.Range("J58").Value = 6100
$ws.Range("K58").Value = 650.5714
$ws.Range("L58").Value = 6100
$ws.Range("M58").Value = -447.5714
$ws.Range("N58").Value = -6506

$ws.Range("H132").Value = 10511.333
$ws.Range("I132").Value = 11304.833
$ws.Range("K132").Value = 33914.499
$ws.Range("M132").Value = -31384.499

$ws.Range("H134").Value = 2465.158
$ws.Range("I134").Value = 2167
$ws.Range("K134").Value = 6501
$ws.Range("M134").Value = -3966

$ws.Range("H136").Value = 2084.6316
$ws.Range("I136").Value = 650.5714
$ws.Range("J136").Value = 6100
$ws.Range("K136").Value = 1951.7142
$ws.Range("L136").Value = 18300
$ws.Range("M136").Value = 598.2857999999999
$ws.Range("N136").Value = -23400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1433.9166
$ws.Range("I44").Value = 332.22223
$ws.Range("J44").Value = 2094.9333
$ws.Range("K44").Value = 996.66669
$ws.Range("L44").Value = 6284.7999
$ws.Range("M44").Value = -598.66669
$ws.Range("N44").Value = -7080.7999

$ws.Range("H69").Value = 2130
$ws.Range("I69").Value = 390
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 1170
$ws.Range("L69").Value = 9000
$ws.Range("M69").Value = -359
$ws.Range("N69").Value = -10622

$ws.Range("H72").Value = 2130
$ws.Range("I72").Value = 390
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 3510
$ws.Range("L72").Value = 27000
$ws.Range("M72").Value = 546
$ws.Range("N72").Value = -35112

$ws.Range("H74").Value = 8146.231
$ws.Range("I74").Value = 4987.5713
$ws.Range("J74").Value = 11831.333
$ws.Range("K74").Value = 14962.7139
$ws.Range("L74").Value = 35493.999
$ws.Range("M74").Value = -13901.7139
$ws.Range("N74").Value = -37615.999

$ws.Range("H77").Value = 8146.231
$ws.Range("I77").Value = 4987.5713
$ws.Range("J77").Value = 11831.333
$ws.Range("K77").Value = 44888.14169999999
$ws.Range("L77").Value = 106481.997
$ws.Range("M77").Value = -39584.14169999999
$ws.Range("N77").Value = -117089.997

$ws.Range("H80").Value = 3100
$ws.Range("J80").Value = 3100
$ws.Range("L80").Value = 9300
$ws.Range("N80").Value = -11172

$ws.Range("H83").Value = 3100
$ws.Range("J83").Value = 3100
$ws.Range("L83").Value = 27900
$ws.Range("N83").Value = -37260

$ws.Range("H87").Value = 13037.375
$ws.Range("I87").Value = 5066.3335
$ws.Range("K87").Value = 15199.0005
$ws.Range("M87").Value = -13951.0005

$ws.Range("H90").Value = 13037.375
$ws.Range("I90").Value = 5066.3335
$ws.Range("K90").Value = 45597.0015
$ws.Range("M90").Value = -39357.0015

$ws.Range("H107").Value = 729.89655
$ws.Range("I107").Value = 219
$ws.Range("J107").Value = 892.4545000000001
$ws.Range("K107").Value = 657
$ws.Range("L107").Value = 2677.3635
$ws.Range("M107").Value = 1263
$ws.Range("N107").Value = -6517.3635

$ws.Range("H133").Value = 6577.778
$ws.Range("J133").Value = 6800
$ws.Range("L133").Value = 20400
$ws.Range("N133").Value = -30520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6060.227
$ws.Range("I132").Value = 5024.9165
$ws.Range("J132").Value = 7302.6
$ws.Range("K132").Value = 15074.7495
$ws.Range("L132").Value = 21907.8
$ws.Range("M132").Value = -12544.7495
$ws.Range("N132").Value = -26967.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1317.2142
$ws.Range("I46").Value = 1213.6666
$ws.Range("J46").Value = 1345.4546
$ws.Range("K46").Value = 1213.6666
$ws.Range("L46").Value = 1345.4546
$ws.Range("M46").Value = -1025.6666
$ws.Range("N46").Value = -1721.4546

$ws.Range("H132").Value = 13522085
$ws.Range("I132").Value = 5524
$ws.Range("J132").Value = 25011162
$ws.Range("K132").Value = 16572
$ws.Range("L132").Value = 75033486
$ws.Range("M132").Value = -14042
$ws.Range("N132").Value = -75038546

$ws.Range("H136").Value = 9707.200000000001
$ws.Range("I136").Value = 2259.3
$ws.Range("J136").Value = 24603
$ws.Range("K136").Value = 6777.900000000001
$ws.Range("L136").Value = 73809
$ws.Range("M136").Value = -4227.900000000001
$ws.Range("N136").Value = -78909

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12248.083
$ws.Range("J62").Value = 16425.285
$ws.Range("L62").Value = 16425.285
$ws.Range("N62").Value = -17673.285

$ws.Range("H64").Value = 24173.854
$ws.Range("J64").Value = 24173.854
$ws.Range("L64").Value = 24173.854
$ws.Range("N64").Value = -24669.854

$ws.Range("H65").Value = 12248.083
$ws.Range("J65").Value = 16425.285
$ws.Range("L65").Value = 82126.425
$ws.Range("N65").Value = -88366.425

$ws.Range("H67").Value = 24173.854
$ws.Range("J67").Value = 24173.854
$ws.Range("L67").Value = 24173.854
$ws.Range("N67").Value = -25889.854

$ws.Range("H132").Value = 3240.0667
$ws.Range("I132").Value = 2410.8333
$ws.Range("J132").Value = 3792.889
$ws.Range("K132").Value = 7232.499899999999
$ws.Range("L132").Value = 11378.667
$ws.Range("M132").Value = -4702.499899999999
$ws.Range("N132").Value = -16438.667

$ws.Range("H136").Value = 1679.875
$ws.Range("I136").Value = 1725.7273
$ws.Range("J136").Value = 1579
$ws.Range("K136").Value = 5177.1819
$ws.Range("L136").Value = 4737
$ws.Range("M136").Value = -2627.1819
$ws.Range("N136").Value = -9837
